$wb = $excel.ActiveWorkbook

# Rename sheets with updated timestamped names
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477893663982"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778961262655"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778961272342"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778961742454"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778962382307"

# Sheet 1 - GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778936249807.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778936480505.csv"
$ws1.Range("B4").Value = "go_stims-16504778936489804.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778936630166.csv"

# Sheet 2 - NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_0-16504778937129874.csv"
$ws2.Range("B3").Value = "ZB-match_1-16504778938690202.csv"
$ws2.Range("B4").Value = "OB-16504778945000145.csv"
$ws2.Range("B5").Value = "TB-1650477896101233.csv"
$ws2.Range("B6").Value = "TB-16504778957020154.csv"
$ws2.Range("B7").Value = "ZB-match_0-16504778937819839.csv"
$ws2.Range("B8").Value = "OB-16504778950609822.csv"
$ws2.Range("B9").Value = "TB-1650477895336981.csv"
$ws2.Range("B10").Value = "OB-16504778946819825.csv"

# Sheet 4 - TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778961422348.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778961302338.csv"
$ws4.Range("B4").Value = "MM_stims-1650477896158263.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778961432347.csv"
$ws4.Range("B6").Value = "MM_stims-16504778961742454.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778961592295.csv"

# Sheet 5 - vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650477896222233.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778961782432.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778962062318.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778961902337.csv"
